$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Cells.Item(2, 4)
$dCell.Value = "'29.106.31"
$dCell.Style = "Normal"
$ws.Range("E2").Value = "  -1.14%  "

$dCell = $ws.Cells.Item(3, 4)
$dCell.Value = "'1.972.81"
$dCell.Style = "Normal"
$ws.Range("E3").Value = "  -0.74%  "

$dCell = $ws.Cells.Item(4, 4)
$dCell.Value = "'1.013"
$dCell.Style = "Normal"
$ws.Range("E4").Value = "  +0.64%  "

$dCell = $ws.Cells.Item(5, 4)
$dCell.Value = "'329.36"
$dCell.Style = "Normal"
$ws.Range("E5").Value = "  -0.12%  "

$dCell = $ws.Cells.Item(6, 4)
$dCell.Value = "'1.012"
$dCell.Style = "Normal"
$ws.Range("E6").Value = "  +0.67%  "

$dCell = $ws.Cells.Item(7, 4)
$dCell.Value = "'0.4958"
$dCell.Style = "Normal"
$ws.Range("E7").Value = "  -0.22%  "

$dCell = $ws.Cells.Item(8, 4)
$dCell.Value = "'0.4203"
$dCell.Style = "Normal"
$ws.Range("E8").Value = "  -0.04%  "

$dCell = $ws.Cells.Item(9, 4)
$dCell.Value = "'54.21"
$dCell.Style = "Normal"
$ws.Range("E9").Value = "  +4.29%  "

$dCell = $ws.Cells.Item(10, 4)
$dCell.Value = "'0.09348"
$dCell.Style = "Normal"
$ws.Range("E10").Value = "  +5.27%  "

$ws.Range("E11").Value = "  -1.82%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$dCell = $ws.Cells.Item(12, 4)
$dCell.Value = "'2.057.22"
$dCell.Style = "Normal"
$ws.Range("E12").Value = "  +4.95%  "

$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$dCell = $ws.Cells.Item(13, 4)
$dCell.Value = "'22.72"
$dCell.Style = "Normal"
$ws.Range("E13").Value = "  -2.60%  "

$dCell = $ws.Cells.Item(14, 4)
$dCell.Value = "'7.877"
$dCell.Style = "Normal"
$ws.Range("E14").Value = "  -2.09%  "

$dCell = $ws.Cells.Item(15, 4)
$dCell.Value = "'6.452"
$dCell.Style = "Normal"
$ws.Range("E15").Value = "  -0.65%  "

$ws.Range("E16").Value = "  +0.75%  "

$dCell = $ws.Cells.Item(17, 4)
$dCell.Value = "'0.00001111"
$dCell.Style = "Normal"
$ws.Range("E17").Value = "  +0.62%  "

$dCell = $ws.Cells.Item(18, 4)
$dCell.Value = "'91.75"
$dCell.Style = "Normal"
$ws.Range("E18").Value = "  -4.44%  "

$dCell = $ws.Cells.Item(19, 4)
$dCell.Value = "'0.06725"
$dCell.Style = "Normal"
$ws.Range("E19").Value = "  +1.56%  "

$dCell = $ws.Cells.Item(20, 4)
$dCell.Value = "'19.16"
$dCell.Style = "Normal"
$ws.Range("E20").Value = "  -2.70%  "

$ws.Range("E21").Value = "  +0.55%  "

$dCell = $ws.Cells.Item(22, 4)
$dCell.Value = "'5.950"
$dCell.Style = "Normal"
$ws.Range("E22").Value = "  +0.01%  "

$dCell = $ws.Cells.Item(23, 4)
$dCell.Value = "'29.130.24"
$dCell.Style = "Normal"
$ws.Range("E23").Value = "  -1.05%  "

$dCell = $ws.Cells.Item(24, 4)
$dCell.Value = "'11.96"
$dCell.Style = "Normal"
$ws.Range("E24").Value = "  +1.02%  "

$dCell = $ws.Cells.Item(25, 4)
$dCell.Value = "'2.266"
$dCell.Style = "Normal"
$ws.Range("E25").Value = "  -0.59%  "

$dCell = $ws.Cells.Item(26, 4)
$dCell.Value = "'2.218.67"
$dCell.Style = "Normal"
$ws.Range("E26").Value = "  +1.30%  "

$dCell = $ws.Cells.Item(27, 4)
$dCell.Value = "'20.76"
$dCell.Style = "Normal"
$ws.Range("E27").Value = "  +1.06%  "

$dCell = $ws.Cells.Item(28, 4)
$dCell.Value = "'156.78"
$dCell.Style = "Normal"
$ws.Range("E28").Value = "  -0.31%  "

$dCell = $ws.Cells.Item(29, 4)
$dCell.Value = "'6.251"
$dCell.Style = "Normal"
$ws.Range("E29").Value = "  -4.08%  "

$dCell = $ws.Cells.Item(30, 4)
$dCell.Value = "'2.263"
$dCell.Style = "Normal"
$ws.Range("E30").Value = "  -3.00%  "

$dCell = $ws.Cells.Item(31, 4)
$dCell.Value = "'127.16"
$dCell.Style = "Normal"
$ws.Range("E31").Value = "  -0.49%  "

$dCell = $ws.Cells.Item(32, 4)
$dCell.Value = "'1.044"
$dCell.Style = "Normal"
$ws.Range("E32").Value = "  -0.59%  "

$dCell = $ws.Cells.Item(33, 4)
$dCell.Value = "'0.09825"
$dCell.Style = "Normal"
$ws.Range("E33").Value = "  -1.02%  "

$dCell = $ws.Cells.Item(34, 4)
$dCell.Value = "'1.499"
$dCell.Style = "Normal"
$ws.Range("E34").Value = "  -4.26%  "

$dCell = $ws.Cells.Item(35, 4)
$dCell.Value = "'5.806"
$dCell.Style = "Normal"
$ws.Range("E35").Value = "  -0.44%  "

$dCell = $ws.Cells.Item(36, 4)
$dCell.Value = "'3.748"
$dCell.Style = "Normal"
$ws.Range("E36").Value = "  -1.04%  "

$dCell = $ws.Cells.Item(37, 4)
$dCell.Value = "'0.02416"
$dCell.Style = "Normal"
$ws.Range("E37").Value = "  -1.28%  "

$dCell = $ws.Cells.Item(38, 4)
$dCell.Value = "'1.326"
$dCell.Style = "Normal"
$ws.Range("E38").Value = "  +3.27%  "

$dCell = $ws.Cells.Item(39, 4)
$dCell.Value = "'0.06411"
$dCell.Style = "Normal"
$ws.Range("E39").Value = "  +1.03%  "

$dCell = $ws.Cells.Item(40, 4)
$dCell.Value = "'9.037"
$dCell.Style = "Normal"
$ws.Range("E40").Value = "  -5.48%  "

$ws.Range("E41").Value = "  -0.46%  "

$dCell = $ws.Cells.Item(42, 4)
$dCell.Value = "'11.49"
$dCell.Style = "Normal"
$ws.Range("E42").Value = "  -2.24%  "

$dCell = $ws.Cells.Item(43, 4)
$dCell.Value = "'0.2001"
$dCell.Style = "Normal"
$ws.Range("E43").Value = "  -3.01%  "

$ws.Range("E44").Value = "  +0.54%  "

$dCell = $ws.Cells.Item(45, 4)
$dCell.Value = "'0.6200"
$dCell.Style = "Normal"
$ws.Range("E45").Value = "  -2.23%  "

$dCell = $ws.Cells.Item(46, 4)
$dCell.Value = "'1.357"
$dCell.Style = "Normal"
$ws.Range("E46").Value = "  +7.14%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$dCell = $ws.Cells.Item(47, 4)
$dCell.Value = "'13.24"
$dCell.Style = "Normal"
$ws.Range("E47").Value = "  -1.30%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$dCell = $ws.Cells.Item(48, 4)
$dCell.Value = "'2.176"
$dCell.Style = "Normal"
$ws.Range("E48").Value = "  -1.76%  "

$dCell = $ws.Cells.Item(49, 4)
$dCell.Value = "'3.488"
$dCell.Style = "Normal"
$ws.Range("E49").Value = "  -1.27%  "

$dCell = $ws.Cells.Item(50, 4)
$dCell.Value = "'0.00000000331"
$dCell.Style = "Normal"
$ws.Range("E50").Value = "  -0.20%  "

$dCell = $ws.Cells.Item(51, 4)
$dCell.Value = "'0.06966"
$dCell.Style = "Normal"
$ws.Range("E51").Value = "  -0.29%  "
